# Misc updates to code
#
# The "Methow River Fawn 02" reach row is removed from the habitat-quality
# table (row 17). All subsequent reach rows shift up by one, which also
# shrinks the sheet's used range from A1:W30 to A1:W29.
#
# In addition, the "Methow River Rattlesnake 06" reach (now row 20 after
# the shift) gets its Flow-SummerBaseFlow_score (column L) corrected from
# 5 to 1, with the dependent HQ_Sum (T) and HQ_Pct (U) totals recalculated
# to match (45 -> 41, 1 -> 0.9111111111111111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Methow River Fawn 02" row entirely; Excel shifts rows 18-30
# up to fill the gap (new row 29 becomes the old row 30 / White River data).
$ws.Rows(17).Delete()

# Correct the Rattlesnake 06 reach's Flow-SummerBaseFlow_score and its
# dependent totals (this reach now lives at row 20 post-shift).
$ws.Range("L20").Value = 1
$ws.Range("T20").Value = 41
$ws.Range("U20").Value = 0.9111111111111111
